# Insert a new data row at row 474 (shifts existing rows 474:495 down to 475:496)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(474).Insert()

# Populate the newly inserted row with its values
$ws.Range("A474").Value = 9
$ws.Range("B474").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C474").Value = "Metropolitana"
$ws.Range("D474").Value = 44939
$ws.Range("E474").Value = 13
$ws.Range("F474").Value = 100112044
$ws.Range("G474").Value = "Perejil"
$ws.Range("H474").Value = "Sin especificar"
$ws.Range("I474").Value = "Primera"
$ws.Range("J474").Value = 70
$ws.Range("K474").Value = 11000
$ws.Range("L474").Value = 13000
$ws.Range("M474").Value = 12000
$ws.Range("N474").Value = "$/docena de atados"
$ws.Range("O474").Value = "Región Metropolitana"
$ws.Range("P474").Value = 4000
$ws.Range("Q474").Value = 3
$ws.Range("R474").Value = "Hortaliza"
